$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# The "形態" (data type) column for the CreateDate / LastUpdate fields is being
# changed from "DATE" to "TIMESTAMP".
$ws.Range("D15").Value = "TIMESTAMP"
$ws.Range("D17").Value = "TIMESTAMP"

# Leave the selection on the last cell that was edited.
$ws.Range("D17").Select() | Out-Null
